# Update the yearly report: drop the "1396/12" twelve-month period column,
# shift all subsequent period data one column to the left (E<-F<-G<-H<-I),
# and append a new "1401/12" twelve-month period with its figures in column I.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header rows (period labels) ---------------------------------------
# Row 8 and Row 24 both hold the same five period headers in E:I.
$headerRows = @(8, 24)
foreach ($r in $headerRows) {
    $ws.Range("E$r").Value = "دوازده ماهه منتهی به 1397/12"
    $ws.Range("F$r").Value = "دوازده ماهه منتهی به 1398/12"
    $ws.Range("G$r").Value = "دوازده ماهه منتهی به 1399/12"
    $ws.Range("H$r").Value = "دوازده ماهه منتهی به 1400/12"
    $ws.Range("I$r").Value = "دوازده ماهه منتهی به 1401/12"
}

# --- Data rows: "هزینه های عمومی و اداری" (general & admin expenses) ---
# Each row's values shift left by one period and a new figure is appended
# in column I for the new 1401/12 period.
$ws.Range("E10").Value = 0
$ws.Range("F10").Value = 0
$ws.Range("G10").Value = 0
$ws.Range("H10").Value = 0
$ws.Range("I10").Value = 0

$ws.Range("E11").Value = 0
$ws.Range("F11").Value = 0
$ws.Range("G11").Value = 0
$ws.Range("H11").Value = 0
$ws.Range("I11").Value = 0

$ws.Range("E12").Value = 0
$ws.Range("F12").Value = 0
$ws.Range("G12").Value = 0
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 0

$ws.Range("E13").Value = 0
$ws.Range("F13").Value = 0
$ws.Range("G13").Value = 0
$ws.Range("H13").Value = 0
$ws.Range("I13").Value = 0

$ws.Range("E14").Value = 0
$ws.Range("F14").Value = 0
$ws.Range("G14").Value = 0
$ws.Range("H14").Value = 0
$ws.Range("I14").Value = 0

$ws.Range("E15").Value = 0
$ws.Range("F15").Value = 0
$ws.Range("G15").Value = 0
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 0

$ws.Range("E16").Value = 0
$ws.Range("F16").Value = 0
$ws.Range("G16").Value = 0
$ws.Range("H16").Value = 459
$ws.Range("I16").Value = 410

$ws.Range("E17").Value = 819
$ws.Range("F17").Value = 12608
$ws.Range("G17").Value = 24453
$ws.Range("H17").Value = 32121
$ws.Range("I17").Value = 48905

$ws.Range("E18").Value = 0
$ws.Range("F18").Value = 0
$ws.Range("G18").Value = 0
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 0

$ws.Range("E19").Value = 4134
$ws.Range("F19").Value = 6675
$ws.Range("G19").Value = 8568
$ws.Range("H19").Value = 28132
$ws.Range("I19").Value = 46038

$ws.Range("E20").Value = 4953
$ws.Range("F20").Value = 19283
$ws.Range("G20").Value = 33021
$ws.Range("H20").Value = 60712
$ws.Range("I20").Value = 95353

# --- Data rows: personnel headcount -------------------------------------
$ws.Range("E26").Value = 15
$ws.Range("F26").Value = 25
$ws.Range("G26").Value = 30
$ws.Range("H26").Value = 35
$ws.Range("I26").Value = 35

$ws.Range("E27").Value = 49
$ws.Range("F27").Value = 60
$ws.Range("G27").Value = 85
$ws.Range("H27").Value = 125
$ws.Range("I27").Value = 145
